{"js": "// Replace the italic title line (\"From Union Square to Rome, Introduction\n// ====...\") and the bold byline (\"By Dorothy Day\") paragraphs with a single,\n// unformatted pandoc-style title-block author line: \"% Dorothy Day\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length < 2) {\n  throw new Error(\"Expected at least two leading paragraphs (title + byline).\");\n}\n\nconst titlePara = paragraphs.items[0];\nconst bylinePara = paragraphs.items[1];\n\n// Drop the byline paragraph (\"By Dorothy Day\") entirely.\nbylinePara.delete();\n\n// Replace the title paragraph's contents/formatting with a plain run that\n// just reads \"% Dorothy Day\" (no italic/bold run properties survive).\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:t xml:space=\"preserve\">% Dorothy Day</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ntitlePara.insertOoxml(ooxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Replace the italic title line (\"From Union Square to Rome, Introduction\n# ====...\") and the bold byline (\"By Dorothy Day\") paragraphs with a single,\n# unformatted pandoc-style title-block author line: \"% Dorothy Day\".\n$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -lt 2) {\n    throw \"Expected at least two leading paragraphs (title + byline).\"\n}\n\n$titlePara  = $d.Paragraphs.Item(1)\n$bylinePara = $d.Paragraphs.Item(2)\n\n# Drop the byline paragraph (\"By Dorothy Day\") entirely, merging its\n# paragraph mark away so the title paragraph becomes the lone survivor.\n$bylinePara.Range.Delete()\n\n# Replace the (now sole) title paragraph's contents/formatting with a plain\n# run that just reads \"% Dorothy Day\" (no italic/bold run properties survive).\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">% Dorothy Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item(1).Range.InsertXML($ooxml)\n"}
